$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($r = 3; $r -le 22; $r++) {
    $hCell = $ws.Cells.Item($r, 8)
    $hCell.Value2 = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($r, 9)
    $iCell.NumberFormat = "@"
    $iCell.Value2 = "04-Nov-2025"
}
